# Update "合肥-漫展信息.xlsx" to reflect refreshed counts scraped at 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both contain the same
# rows of data (the former is a subset/filtered view, the latter the full
# list), so the same cell edits are applied to both worksheets.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("C3").Value = "合肥·首届AS运动番Only（取消）"
$ws1.Range("F2").Value = 127
$ws1.Range("F3").Value = 234
$ws1.Range("F5").Value = 6659
$ws1.Range("F9").Value = 6125
$ws1.Range("F12").Value = 1245
$ws1.Range("F13").Value = 1245
$ws1.Range("F16").Value = 389
$ws1.Range("F22").Value = 4445
$ws1.Range("F23").Value = 52
$ws1.Range("F24").Value = 24
$ws1.Range("F26").Value = 45

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("C3").Value = "合肥·首届AS运动番Only（取消）"
$ws4.Range("F2").Value = 127
$ws4.Range("F3").Value = 234
$ws4.Range("F5").Value = 6659
$ws4.Range("F9").Value = 6125
$ws4.Range("F12").Value = 1245
$ws4.Range("F13").Value = 1245
$ws4.Range("F16").Value = 389
$ws4.Range("F22").Value = 4445
$ws4.Range("F24").Value = 52
$ws4.Range("F25").Value = 24
$ws4.Range("F27").Value = 45
